# Natmi following Dr Hou advice
# Update the LR-pair results table: recompute rows for the full 3x3
# Sending-cluster x Target-cluster matrix (ECs/FAPs/sCs), replacing the
# previous (partial) 3x2 matrix of rows 2-7 and adding new rows 8-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @{
  2 = @("ECs", "L1cam", "Alcam", "ECs", 3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 950.7569416343759, 8556.812474709384, 0.7230173910559203, 0.7230173910559203)
  3 = @("ECs", "L1cam", "Alcam", "FAPs", 3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 10.78136413674756, 97.032277230728, 0.008198850230612124, 0.008198850230612124)
  4 = @("ECs", "L1cam", "Alcam", "sCs", 3, 1, 22.59487733333333, 67.784632, 0.7395019553569895, 0.7395019553569895, 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 10.89558877329333, 98.06029895963999, 0.008285714070457112, 0.008285714070457112)
  5 = @("FAPs", "L1cam", "Alcam", "ECs", 3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 12.817033368828, 115.353300319452, 0.009746905461953911, 0.009746905461953911)
  6 = @("FAPs", "L1cam", "Alcam", "FAPs", 3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 0.1453421982537778, 1.308079784284, 0.0001105276568490102, 0.0001105276568490102)
  7 = @("FAPs", "L1cam", "Alcam", "sCs", 3, 1, 0.3045986666666667, 0.9137960000000001, 0.00996913177602551, 0.00996913177602551, 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 0.1468820460466667, 1.32193841442, 0.0001116986572225903, 0.0001116986572225903)
  8 = @("sCs", "L1cam", "Alcam", "ECs", 3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 3, 1, 42.078429, 126.235287, 0.9777085588730982, 0.9777085588730982, 322.098002936874, 2898.882026431866, 0.2449442623552238, 0.2449442623552238)
  9 = @("sCs", "L1cam", "Alcam", "FAPs", 3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 2, 0.6666666666666666, 0.4771596666666666, 1.431479, 0.01108698925164327, 0.01108698925164327, 3.652516963391333, 32.872652670522, 0.002777611364182135, 0.002777611364182135)
  10 = @("sCs", "L1cam", "Alcam", "sCs", 3, 1, 7.654706, 22.964118, 0.2505289128669849, 0.2505289128669849, 3, 1, 0.482215, 1.446645, 0.01120445187525872, 0.01120445187525872, 3.69121405379, 33.22092648411, 0.002807039147579018, 0.002807039147579018)
}

foreach ($r in ($data.Keys | Sort-Object)) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $rowVals[$i]
    }
}
